$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 187.75
$ws.Cells.Item(39, 9).Value = 117.5
$ws.Cells.Item(39, 11).Value = 352.5
$ws.Cells.Item(39, 13).Value = -56.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 5085.4
$ws.Cells.Item(106, 9).Value = 4787
$ws.Cells.Item(106, 11).Value = 4787
$ws.Cells.Item(106, 13).Value = -4156

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 84147
$ws.Cells.Item(107, 9).Value = 125620.25
$ws.Cells.Item(107, 10).Value = 1200.5
$ws.Cells.Item(107, 11).Value = 125620.25
$ws.Cells.Item(107, 12).Value = 1200.5
$ws.Cells.Item(107, 13).Value = -123700.25
$ws.Cells.Item(107, 14).Value = -5040.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1286.375
$ws.Cells.Item(111, 9).Value = 816.8
$ws.Cells.Item(111, 10).Value = 2069
$ws.Cells.Item(111, 11).Value = 2450.4
$ws.Cells.Item(111, 12).Value = 6207
$ws.Cells.Item(111, 13).Value = 616.6000000000004
$ws.Cells.Item(111, 14).Value = -12341

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 4478.3335
$ws.Cells.Item(113, 9).Value = 3866.5
$ws.Cells.Item(113, 10).Value = 5702
$ws.Cells.Item(113, 11).Value = 3866.5
$ws.Cells.Item(113, 12).Value = 5702
$ws.Cells.Item(113, 13).Value = -612.5
$ws.Cells.Item(113, 14).Value = -12210

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2661.2273
$ws.Cells.Item(138, 9).Value = 1828.75
$ws.Cells.Item(138, 10).Value = 3274.6316
$ws.Cells.Item(138, 11).Value = 5486.25
$ws.Cells.Item(138, 12).Value = 9823.8948
$ws.Cells.Item(138, 13).Value = -346.25
$ws.Cells.Item(138, 14).Value = -20103.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9215.596
$ws.Cells.Item(32, 9).Value = 5624.2856
$ws.Cells.Item(32, 11).Value = 5624.2856
$ws.Cells.Item(32, 13).Value = -5337.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 30846.314
$ws.Cells.Item(74, 9).Value = 32847.094
$ws.Cells.Item(74, 11).Value = 32847.094
$ws.Cells.Item(74, 13).Value = -31973.094

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 30846.314
$ws.Cells.Item(77, 9).Value = 32847.094
$ws.Cells.Item(77, 11).Value = 164235.47
$ws.Cells.Item(77, 13).Value = -159867.47

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(105, 8).Value = 118999.5
$ws.Cells.Item(105, 10).Value = 118999.5
$ws.Cells.Item(105, 12).Value = 118999.5
$ws.Cells.Item(105, 14).Value = -125987.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3382.1538
$ws.Cells.Item(132, 9).Value = 3024.2222
$ws.Cells.Item(132, 11).Value = 9072.6666
$ws.Cells.Item(132, 13).Value = -6542.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(140, 8).Value = 162797
$ws.Cells.Item(140, 10).Value = 162797
$ws.Cells.Item(140, 12).Value = 162797
$ws.Cells.Item(140, 14).Value = -173157

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 105024.65
$ws.Cells.Item(20, 9).Value = 148735.88
$ws.Cells.Item(20, 11).Value = 148735.88
$ws.Cells.Item(20, 13).Value = -148488.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 643.43475
$ws.Cells.Item(94, 9).Value = 401.1613
$ws.Cells.Item(94, 11).Value = 401.1613
$ws.Cells.Item(94, 13).Value = 49.83870000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1877563
$ws.Cells.Item(99, 9).Value = 2713.0625
$ws.Cells.Item(99, 11).Value = 2713.0625
$ws.Cells.Item(99, 13).Value = -1215.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4402.067
$ws.Cells.Item(31, 9).Value = 2849.25
$ws.Cells.Item(31, 10).Value = 5644.32
$ws.Cells.Item(31, 11).Value = 2849.25
$ws.Cells.Item(31, 12).Value = 5644.32
$ws.Cells.Item(31, 13).Value = -2554.25
$ws.Cells.Item(31, 14).Value = -6234.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4402.067
$ws.Cells.Item(34, 9).Value = 2849.25
$ws.Cells.Item(34, 10).Value = 5644.32
$ws.Cells.Item(34, 11).Value = 2849.25
$ws.Cells.Item(34, 12).Value = 5644.32
$ws.Cells.Item(34, 13).Value = -2647.25
$ws.Cells.Item(34, 14).Value = -6048.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 333335970
$ws.Cells.Item(86, 9).Value = 333335970
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 333335970
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).Value = -333334847

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 333335970
$ws.Cells.Item(89, 9).Value = 333335970
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 1666679850
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).Value = -1666674234

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 67549.53
$ws.Cells.Item(105, 9).Value = 125261.555
$ws.Cells.Item(105, 11).Value = 125261.555
$ws.Cells.Item(105, 13).Value = -123514.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2669.8
$ws.Cells.Item(132, 9).Value = 2462.25
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 7386.75
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -4856.75
$ws.Cells.Item(132, 14).Value = -15560

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 126788.5
$ws.Cells.Item(134, 9).Value = 1361.8
$ws.Cells.Item(134, 10).Value = 335833
$ws.Cells.Item(134, 11).Value = 4085.4
$ws.Cells.Item(134, 12).Value = 1007499
$ws.Cells.Item(134, 13).Value = -1550.4
$ws.Cells.Item(134, 14).Value = -1012569

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 181.18518
$ws.Cells.Item(2, 9).Value = 144.63637
$ws.Cells.Item(2, 11).Value = 867.81822
$ws.Cells.Item(2, 13).Value = -754.81822

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 5240.9614
$ws.Cells.Item(112, 9).Value = 2488.6667
$ws.Cells.Item(112, 10).Value = 5599.9565
$ws.Cells.Item(112, 11).Value = 7466.000100000001
$ws.Cells.Item(112, 12).Value = 16799.8695
$ws.Cells.Item(112, 13).Value = -6358.000100000001
$ws.Cells.Item(112, 14).Value = -19015.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1269
$ws.Cells.Item(113, 9).Value = 1132
$ws.Cells.Item(113, 10).Value = 1303.25
$ws.Cells.Item(113, 11).Value = 3396
$ws.Cells.Item(113, 12).Value = 3909.75
$ws.Cells.Item(113, 13).Value = -1226
$ws.Cells.Item(113, 14).Value = -8249.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 36605.723
$ws.Cells.Item(131, 10).Value = 2425.0476
$ws.Cells.Item(131, 12).Value = 7275.1428
$ws.Cells.Item(131, 14).Value = -17355.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6666.263
$ws.Cells.Item(70, 9).Value = 6205.077
$ws.Cells.Item(70, 11).Value = 6205.077
$ws.Cells.Item(70, 13).Value = -5935.077

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 6666.263
$ws.Cells.Item(73, 9).Value = 6205.077
$ws.Cells.Item(73, 11).Value = 6205.077
$ws.Cells.Item(73, 13).Value = -5269.077

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 100007640
$ws.Cells.Item(80, 9).Value = 200002750
$ws.Cells.Item(80, 11).Value = 200002750
$ws.Cells.Item(80, 13).Value = -200001752

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 100007640
$ws.Cells.Item(83, 9).Value = 200002750
$ws.Cells.Item(83, 11).Value = 1000013750
$ws.Cells.Item(83, 13).Value = -1000008758

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 6561.5
$ws.Cells.Item(97, 9).Value = 7082
$ws.Cells.Item(97, 10).Value = 5000
$ws.Cells.Item(97, 11).Value = 7082
$ws.Cells.Item(97, 12).Value = 5000
$ws.Cells.Item(97, 13).Value = -6586
$ws.Cells.Item(97, 14).Value = -5992

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1258.3334
$ws.Cells.Item(102, 9).Value = 1283.04
$ws.Cells.Item(102, 11).Value = 1283.04
$ws.Cells.Item(102, 13).Value = 338.96

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3727.0908
$ws.Cells.Item(126, 9).Value = 2899.3333
$ws.Cells.Item(126, 11).Value = 8697.999899999999
$ws.Cells.Item(126, 13).Value = -6227.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2956.5454
$ws.Cells.Item(132, 9).Value = 2289
$ws.Cells.Item(132, 11).Value = 6867
$ws.Cells.Item(132, 13).Value = -4337

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 14).Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 304357.16
$ws.Cells.Item(68, 9).Value = 304357.16
$ws.Cells.Item(68, 11).Value = 304357.16
$ws.Cells.Item(68, 13).Value = -303608.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 304357.16
$ws.Cells.Item(71, 9).Value = 304357.16
$ws.Cells.Item(71, 11).Value = 1521785.8
$ws.Cells.Item(71, 13).Value = -1518041.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1289.1111
$ws.Cells.Item(82, 9).Value = 664.3570999999999
$ws.Cells.Item(82, 11).Value = 664.3570999999999
$ws.Cells.Item(82, 13).Value = -303.3570999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1289.1111
$ws.Cells.Item(85, 9).Value = 664.3570999999999
$ws.Cells.Item(85, 11).Value = 664.3570999999999
$ws.Cells.Item(85, 13).Value = 583.6429000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 8999.5
$ws.Cells.Item(93, 9).Value = 10000
$ws.Cells.Item(93, 11).Value = 10000
$ws.Cells.Item(93, 13).Value = -8752

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2500
$ws.Cells.Item(132, 9).Value = 2500
$ws.Cells.Item(132, 11).Value = 7500
$ws.Cells.Item(132, 13).Value = -4970

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 6101.6763
$ws.Cells.Item(136, 9).Value = 6247.9
$ws.Cells.Item(136, 10).Value = 5892.7856
$ws.Cells.Item(136, 11).Value = 18743.7
$ws.Cells.Item(136, 12).Value = 17678.3568
$ws.Cells.Item(136, 13).Value = -16193.7
$ws.Cells.Item(136, 14).Value = -22778.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 50749.5
$ws.Cells.Item(81, 9).Value = 1499
$ws.Cells.Item(81, 11).Value = 2998
$ws.Cells.Item(81, 13).Value = -1937

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 50749.5
$ws.Cells.Item(84, 9).Value = 1499
$ws.Cells.Item(84, 11).Value = 14990
$ws.Cells.Item(84, 13).Value = -9686

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1979582.8
$ws.Cells.Item(132, 9).Value = 3304.0667
$ws.Cells.Item(132, 10).Value = 6214465.5
$ws.Cells.Item(132, 11).Value = 9912.2001
$ws.Cells.Item(132, 12).Value = 18643396.5
$ws.Cells.Item(132, 13).Value = -7382.2001
$ws.Cells.Item(132, 14).Value = -18648456.5
